$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.06446533333333333
$ws.Range("H2").Value = 0.193396
$ws.Range("I2").Value = 0.02693738696927793
$ws.Range("J2").Value = 0.02693738696927793
$ws.Range("M2").Value = 1.334383666666667
$ws.Range("N2").Value = 4.003151
$ws.Range("O2").Value = 0.1312069045987744
$ws.Range("P2").Value = 0.1312069045987744
$ws.Range("Q2").Value = 0.08602148786622221
$ws.Range("R2").Value = 0.774193390796
$ws.Range("S2").Value = 0.003534371162218317
$ws.Range("T2").Value = 0.003534371162218317
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.06446533333333333
$ws.Range("H3").Value = 0.193396
$ws.Range("I3").Value = 0.02693738696927793
$ws.Range("J3").Value = 0.02693738696927793
$ws.Range("N3").Value = 7.432386999999999
$ws.Range("O3").Value = 0.2436032245723858
$ws.Range("P3").Value = 0.2436032245723858
$ws.Range("Q3").Value = 0.1597104351391111
$ws.Range("R3").Value = 1.437393916252
$ws.Range("S3").Value = 0.00656203432727027
$ws.Range("T3").Value = 0.006562034327270271
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.06446533333333333
$ws.Range("H4").Value = 0.193396
$ws.Range("I4").Value = 0.02693738696927793
$ws.Range("J4").Value = 0.02693738696927793
$ws.Range("M4").Value = 6.358226000000001
$ws.Range("N4").Value = 19.074678
$ws.Range("O4").Value = 0.6251898708288398
$ws.Range("P4").Value = 0.6251898708288398
$ws.Range("Q4").Value = 0.4098851584986667
$ws.Range("R4").Value = 3.688966426488001
$ws.Range("S4").Value = 0.01684098147978934
$ws.Range("T4").Value = 0.01684098147978934
$ws.Range("I5").Value = 0.7704314695358874
$ws.Range("J5").Value = 0.7704314695358874
$ws.Range("M5").Value = 1.334383666666667
$ws.Range("N5").Value = 4.003151
$ws.Range("O5").Value = 0.1312069045987744
$ws.Range("P5").Value = 0.1312069045987744
$ws.Range("Q5").Value = 2.460285453226111
$ws.Range("R5").Value = 22.142569079035
$ws.Range("S5").Value = 0.1010859283232887
$ws.Range("T5").Value = 0.1010859283232887
$ws.Range("I6").Value = 0.7704314695358874
$ws.Range("J6").Value = 0.7704314695358874
$ws.Range("N6").Value = 7.432386999999999
$ws.Range("O6").Value = 0.2436032245723858
$ws.Range("P6").Value = 0.2436032245723858
$ws.Range("S6").Value = 0.187679590290984
$ws.Range("T6").Value = 0.187679590290984
$ws.Range("I7").Value = 0.7704314695358874
$ws.Range("J7").Value = 0.7704314695358874
$ws.Range("M7").Value = 6.358226000000001
$ws.Range("N7").Value = 19.074678
$ws.Range("O7").Value = 0.6251898708288398
$ws.Range("P7").Value = 0.6251898708288398
$ws.Range("Q7").Value = 11.72305336680334
$ws.Range("R7").Value = 105.50748030123
$ws.Range("S7").Value = 0.4816659509216147
$ws.Range("T7").Value = 0.4816659509216147
$ws.Range("G8").Value = 0.4849276666666666
$ws.Range("H8").Value = 1.454783
$ws.Range("I8").Value = 0.2026311434948347
$ws.Range("J8").Value = 0.2026311434948347
$ws.Range("M8").Value = 1.334383666666667
$ws.Range("N8").Value = 4.003151
$ws.Range("O8").Value = 0.1312069045987744
$ws.Range("P8").Value = 0.1312069045987744
$ws.Range("Q8").Value = 0.6470795579147778
$ws.Range("R8").Value = 5.823716021232999
$ws.Range("S8").Value = 0.02658660511326734
$ws.Range("T8").Value = 0.02658660511326734
$ws.Range("G9").Value = 0.4849276666666666
$ws.Range("H9").Value = 1.454783
$ws.Range("I9").Value = 0.2026311434948347
$ws.Range("J9").Value = 0.2026311434948347
$ws.Range("N9").Value = 7.432386999999999
$ws.Range("O9").Value = 0.2436032245723858
$ws.Range("P9").Value = 0.2436032245723858
$ws.Range("Q9").Value = 1.201390028557889
$ws.Range("R9").Value = 10.812510257021
$ws.Range("S9").Value = 0.04936159995413154
$ws.Range("T9").Value = 0.04936159995413156
$ws.Range("G10").Value = 0.4849276666666666
$ws.Range("H10").Value = 1.454783
$ws.Range("I10").Value = 0.2026311434948347
$ws.Range("J10").Value = 0.2026311434948347
$ws.Range("M10").Value = 6.358226000000001
$ws.Range("N10").Value = 19.074678
$ws.Range("O10").Value = 0.6251898708288398
$ws.Range("P10").Value = 0.6251898708288398
$ws.Range("Q10").Value = 3.083279698319334
$ws.Range("R10").Value = 27.749517284874
$ws.Range("S10").Value = 0.1266829384274358
$ws.Range("T10").Value = 0.1266829384274358
